$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, pushing the existing rows 220-239
# down to 221-240 (this also keeps all their data/formatting intact).
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new weekly record.
$ws.Cells.Item(220, 1).Value = 4
$ws.Cells.Item(220, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(220, 3).Value = "Los Lagos"
$ws.Cells.Item(220, 4).Value = 44449
$ws.Cells.Item(220, 5).Value = 10
$ws.Cells.Item(220, 6).Value = "Fruta"
$ws.Cells.Item(220, 7).Value = 100108
$ws.Cells.Item(220, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(220, 9).Value = 100108006
$ws.Cells.Item(220, 10).Value = "Plátano"
$ws.Cells.Item(220, 11).Value = "Sin especificar"
$ws.Cells.Item(220, 12).Value = "Primera Pintón"
$ws.Cells.Item(220, 13).Value = 800
$ws.Cells.Item(220, 14).Value = 25000
$ws.Cells.Item(220, 15).Value = 25000
$ws.Cells.Item(220, 16).Value = 25000
$ws.Cells.Item(220, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(220, 18).Value = "Ecuador"
$ws.Cells.Item(220, 19).Value = 1250
$ws.Cells.Item(220, 20).Value = 20
